$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Add two new rows to the table (this expands table ref / autoFilter / dimension)
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# ---- Copy number formats (styles) from existing template cells so that ----
# ---- new cells reuse the same style indices instead of creating new ones ----

# Row 30 formats: B,C,F -> datetime style; P -> date style; S,T,V,W,X -> quote-prefixed style
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$ws.Range("P2").Copy() | Out-Null
$ws.Range("P30").PasteSpecial(-4122) | Out-Null
$ws.Range("S2").Copy() | Out-Null
$ws.Range("S30").PasteSpecial(-4122) | Out-Null
$ws.Range("T30").PasteSpecial(-4122) | Out-Null
$ws.Range("V30").PasteSpecial(-4122) | Out-Null
$ws.Range("W30").PasteSpecial(-4122) | Out-Null
$ws.Range("X30").PasteSpecial(-4122) | Out-Null

# Row 31 formats (note: U31, Y31, AE31 stay untouched/absent on purpose)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null
$ws.Range("P2").Copy() | Out-Null
$ws.Range("P31").PasteSpecial(-4122) | Out-Null
$ws.Range("S2").Copy() | Out-Null
$ws.Range("S31").PasteSpecial(-4122) | Out-Null
$ws.Range("T31").PasteSpecial(-4122) | Out-Null
$ws.Range("V31").PasteSpecial(-4122) | Out-Null
$ws.Range("W31").PasteSpecial(-4122) | Out-Null
$ws.Range("X31").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------- Row 30 values ----------------------------
$ws.Cells.Item(30,1).Value  = 31
$ws.Cells.Item(30,2).Value  = 45476.600208333337
$ws.Cells.Item(30,3).Value  = 45476.612546296295
$ws.Cells.Item(30,4).Value  = "anonymous"
$ws.Cells.Item(30,7).Value  = "Si"
$ws.Cells.Item(30,8).Value  = "Administrativo Supernumerario"
$ws.Cells.Item(30,9).Value  = "Calle 72"
$ws.Cells.Item(30,10).Value = "Mujer"
$ws.Cells.Item(30,11).Value = "28 a 40 años"
$ws.Cells.Item(30,12).Value = "Ninguna"
$ws.Cells.Item(30,13).Value = "Sin pertenencia étnica"
$ws.Cells.Item(30,14).Value = "GRUPO DE ASEGURAMIENTO DE LA CALIDAD"
$ws.Cells.Item(30,15).Value = "Angie Manrique "
$ws.Cells.Item(30,16).Value = 45476
$ws.Cells.Item(30,17).Value = "Servicio de transporte local"
$ws.Cells.Item(30,18).Value = "Néstor Raúl Acosta"
$ws.Cells.Item(30,19).Value = "'5"
$ws.Cells.Item(30,20).Value = "'5"
$ws.Cells.Item(30,21).Value = "Excelente "
$ws.Cells.Item(30,22).Value = "'5"
$ws.Cells.Item(30,23).Value = "'5"
$ws.Cells.Item(30,24).Value = "'5"
$ws.Cells.Item(30,25).Value = "Excelente "
$ws.Cells.Item(30,26).Value = "Si"
$ws.Cells.Item(30,27).Value = "Si"
$ws.Cells.Item(30,28).Value = "Si"
$ws.Cells.Item(30,29).Value = "Si"
$ws.Cells.Item(30,30).Value = "Si"
$ws.Cells.Item(30,31).Value = "Ninguna "

# ---------------------------- Row 31 values ----------------------------
$ws.Cells.Item(31,1).Value  = 32
$ws.Cells.Item(31,2).Value  = 45482.464513888888
$ws.Cells.Item(31,3).Value  = 45482.466122685182
$ws.Cells.Item(31,4).Value  = "anonymous"
$ws.Cells.Item(31,7).Value  = "Si"
$ws.Cells.Item(31,8).Value  = "Administrativo Provisional"
$ws.Cells.Item(31,9).Value  = "Calle 79"
$ws.Cells.Item(31,10).Value = "Mujer"
$ws.Cells.Item(31,11).Value = "40 a 60 años"
$ws.Cells.Item(31,12).Value = "Ninguna"
$ws.Cells.Item(31,13).Value = "Sin pertenencia étnica"
$ws.Cells.Item(31,14).Value = "OFICINA DE CONTROL INTERNO"
$ws.Cells.Item(31,15).Value = "Claudia Rojas "
$ws.Cells.Item(31,16).Value = 45469
$ws.Cells.Item(31,17).Value = "Servicio de transporte local"
$ws.Cells.Item(31,18).Value = "Jorge Enrique Mogollón Montañez"
$ws.Cells.Item(31,19).Value = "'4"
$ws.Cells.Item(31,20).Value = "'3"
$ws.Cells.Item(31,22).Value = "'4"
$ws.Cells.Item(31,23).Value = "'4"
$ws.Cells.Item(31,24).Value = "'4"
$ws.Cells.Item(31,26).Value = "Si"
$ws.Cells.Item(31,27).Value = "Si"
$ws.Cells.Item(31,28).Value = "Si"
$ws.Cells.Item(31,29).Value = "No"
$ws.Cells.Item(31,30).Value = "No"
